# Insert a new data row at row 244 (pushing the existing rows 244-337 down
# to 245-338) and populate it with a new "Frutilla" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 244..337 down to 245..338 by inserting a blank row at 244.
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row 244 with the new record.
$ws.Cells.Item(244, 1).Value  = 7
$ws.Cells.Item(244, 2).Value  = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(244, 3).Value  = 'Ñuble'
$ws.Cells.Item(244, 4).Value  = 44837
$ws.Cells.Item(244, 5).Value  = 16
$ws.Cells.Item(244, 6).Value  = 'Fruta'
$ws.Cells.Item(244, 7).Value  = 100101
$ws.Cells.Item(244, 8).Value  = 'Berries'
$ws.Cells.Item(244, 9).Value  = 100112025
$ws.Cells.Item(244, 10).Value = 'Frutilla'
$ws.Cells.Item(244, 11).Value = 'Sin especificar'
$ws.Cells.Item(244, 12).Value = 'Primera'
$ws.Cells.Item(244, 13).Value = 120
$ws.Cells.Item(244, 14).Value = 14000
$ws.Cells.Item(244, 15).Value = 15000
$ws.Cells.Item(244, 16).Value = 14500
$ws.Cells.Item(244, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(244, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(244, 19).Value = 2071
$ws.Cells.Item(244, 20).Value = 7
